# Insert two new data rows at row 34 (this pushes the existing rows 34-111
# down to rows 36-113, matching the target dimension A1:T113).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(34).Insert()
$ws.Rows.Item(34).Insert()

# --- New row 34 ---
$ws.Range("A34").Value = 4
$ws.Range("B34").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C34").Value = "Los Lagos"
$ws.Range("D34").Value2 = 44519
$ws.Range("E34").Value = 10
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100108
$ws.Range("H34").Value = "Tropicales y subtropicales"
$ws.Range("I34").Value = 100108002
$ws.Range("J34").Value = "Mango"
$ws.Range("K34").Value = "Sin especificar"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 200
$ws.Range("N34").Value = 7500
$ws.Range("O34").Value = 8000
$ws.Range("P34").Value = 7750
$ws.Range("Q34").Value = "$/bandeja 4 kilos"
$ws.Range("R34").Value = "Perú"
$ws.Range("S34").Value = 1938
$ws.Range("T34").Value = 4

# --- New row 35 ---
$ws.Range("A35").Value = 4
$ws.Range("B35").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C35").Value = "Los Lagos"
$ws.Range("D35").Value2 = 44519
$ws.Range("E35").Value = 10
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100108
$ws.Range("H35").Value = "Tropicales y subtropicales"
$ws.Range("I35").Value = 100108002
$ws.Range("J35").Value = "Mango"
$ws.Range("K35").Value = "Sin especificar"
$ws.Range("L35").Value = "Segunda"
$ws.Range("M35").Value = 100
$ws.Range("N35").Value = 5000
$ws.Range("O35").Value = 5000
$ws.Range("P35").Value = 5000
$ws.Range("Q35").Value = "$/bandeja 4 kilos"
$ws.Range("R35").Value = "Perú"
$ws.Range("S35").Value = 1250
$ws.Range("T35").Value = 4
